$d = $word.ActiveDocument

# The homework title currently reads "Homework 2: Software Processes".
# Per the commit ("Associated Press-style headlines"), headline-case
# capitalization is dropped for ordinary words after the first, so
# "Processes" becomes "processes". MatchCase ensures we only touch the
# exact-cased title text (there is already an unrelated lowercase
# "processes" elsewhere in the document body that must stay untouched).
$d.Content.Find.Execute("Processes", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "processes", 2)
